# Added new run (run_id=2, "test") and new experiment (experiment_id=2,
# "Test experiment") to the workbook, plus the associated specification
# rows that tie the new run/experiment back to the existing heuristics.

$wb = $excel.ActiveWorkbook

$wsRunDesc  = $wb.Worksheets.Item("run-description")
$wsRunSpec  = $wb.Worksheets.Item("run-specification")
$wsExpDesc  = $wb.Worksheets.Item("experiment-description")
$wsExpSpec  = $wb.Worksheets.Item("experiment-specification")

# --- run-description: new run_id 2, "test" -------------------------------
$wsRunDesc.Range("A3").Value = 2
$wsRunDesc.Range("B3").Value = "test"
$wsRunDesc.Range("C3").Value = "Test run"

# --- run-specification: tie run_id 2 to both experiments -----------------
$wsRunSpec.Range("A3").Value = 2
$wsRunSpec.Range("B3").Value = 1

$wsRunSpec.Range("A4").Value = 2
$wsRunSpec.Range("B4").Value = 2

# --- experiment-description: new experiment_id 2, "Test experiment" ------
$wsExpDesc.Range("A3").Value = 2
$wsExpDesc.Range("B3").Value = "Test experiment"
$wsExpDesc.Range("C3").Value = "extended"

# --- experiment-specification: heuristic + transition/length_of_stay rows
$wsExpSpec.Range("A11").Value = 2
$wsExpSpec.Range("B11").Value = "heuristic"
$wsExpSpec.Range("D11").Value = "heuristic_1"

$wsExpSpec.Range("A12").Value = 2
$wsExpSpec.Range("B12").Value = "transition"
$wsExpSpec.Range("C12").Value = "home-green"
$wsExpSpec.Range("D12").Value = "age_simple"

$wsExpSpec.Range("A13").Value = 2
$wsExpSpec.Range("B13").Value = "transition"
$wsExpSpec.Range("C13").Value = "home-red"
$wsExpSpec.Range("D13").Value = "age_simple"

$wsExpSpec.Range("A14").Value = 2
$wsExpSpec.Range("B14").Value = "length_of_stay"
$wsExpSpec.Range("C14").Value = "home-green"
$wsExpSpec.Range("D14").Value = "age_simple"

$wsExpSpec.Range("A15").Value = 2
$wsExpSpec.Range("B15").Value = "length_of_stay"
$wsExpSpec.Range("C15").Value = "home-red"
$wsExpSpec.Range("D15").Value = "age_simple"

$wsExpSpec.Range("A16").Value = 2
$wsExpSpec.Range("B16").Value = "transition"
$wsExpSpec.Range("C16").Value = "inpatient_ward"
$wsExpSpec.Range("D16").Value = "age_simple"

$wsExpSpec.Range("A17").Value = 2
$wsExpSpec.Range("B17").Value = "length_of_stay"
$wsExpSpec.Range("C17").Value = "inpatient_ward"
$wsExpSpec.Range("D17").Value = "none"

$wsExpSpec.Range("A18").Value = 2
$wsExpSpec.Range("B18").Value = "transition"
$wsExpSpec.Range("C18").Value = "intensive_care_unit"
$wsExpSpec.Range("D18").Value = "age_simple"

$wsExpSpec.Range("A19").Value = 2
$wsExpSpec.Range("B19").Value = "length_of_stay"
$wsExpSpec.Range("C19").Value = "intensive_care_unit"
$wsExpSpec.Range("D19").Value = "none"

# --- selections, matching where the user last clicked on each sheet ------
$wsRunDesc.Range("E15").Select()
$wsRunSpec.Range("C11").Select()
$wsExpSpec.Range("D19").Select()

# experiment-description ends up the active sheet/selection when the file
# was saved.
$wsExpDesc.Activate()
$wsExpDesc.Range("E22").Select()
